$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-5
for ($r = 2; $r -le 4; $r++) {
    $ws.Range("B$r").Value = 3.286832544864788
    $ws.Range("C$r").Value = 1.655778082260271
    $ws.Range("D$r").Value = 0.1494219747398047
    $ws.Range("E$r").Value = 0.4942365360607697
    $ws.Range("F$r").Value = 1
    $ws.Range("G$r").Value = 5.586269137925634
}

$ws.Range("B5").Value = 3.286832544864788
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.189590430959694
